$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value2 = 434.2
$ws.Range("I15").Value2 = 434.2
$ws.Range("K15").Value2 = 1302.6
$ws.Range("M15").Value2 = -1133.6
$ws.Range("H106").Value2 = 74076840
$ws.Range("I106").Value2 = 27780678
$ws.Range("K106").Value2 = 27780678
$ws.Range("M106").Value2 = -27780047
$ws.Range("H107").Value2 = 8333921.5
$ws.Range("I107").Value2 = 13158173
$ws.Range("J107").Value2 = 1123.1818
$ws.Range("K107").Value2 = 13158173
$ws.Range("L107").Value2 = 1123.1818
$ws.Range("M107").Value2 = -13156253
$ws.Range("N107").Value2 = -4963.1818
$ws.Range("H116").Value2 = 6914.364
$ws.Range("I116").Value2 = 9408.462
$ws.Range("K116").Value2 = 9408.462
$ws.Range("M116").Value2 = -5966.462
$ws.Range("H132").Value2 = 6173548.5
$ws.Range("I132").Value2 = 447.30435
$ws.Range("J132").Value2 = 41668880
$ws.Range("K132").Value2 = 1341.91305
$ws.Range("L132").Value2 = 125006640
$ws.Range("M132").Value2 = 1188.08695
$ws.Range("N132").Value2 = -125011700
$ws.Range("H137").Value2 = 1316.6222
$ws.Range("I137").Value2 = 1027.2433
$ws.Range("J137").Value2 = 2655
$ws.Range("K137").Value2 = 3081.7299
$ws.Range("L137").Value2 = 7965
$ws.Range("M137").Value2 = -531.7299000000003
$ws.Range("N137").Value2 = -13065

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1240.909
$ws.Range("I2").Value2 = 801.5714
$ws.Range("K2").Value2 = 801.5714
$ws.Range("M2").Value2 = -688.5714
$ws.Range("H32").Value2 = 5687.3496
$ws.Range("I32").Value2 = 4062.3088
$ws.Range("K32").Value2 = 4062.3088
$ws.Range("M32").Value2 = -3775.3088
$ws.Range("H61").Value2 = 4477.4414
$ws.Range("I61").Value2 = 4636.5483
$ws.Range("J61").Value2 = 2833.3333
$ws.Range("K61").Value2 = 4636.5483
$ws.Range("L61").Value2 = 2833.3333
$ws.Range("M61").Value2 = -4424.5483
$ws.Range("N61").Value2 = -3257.3333
$ws.Range("H74").Value2 = 25001872
$ws.Range("I74").Value2 = 1701.7858
$ws.Range("J74").Value2 = 83335600
$ws.Range("K74").Value2 = 1701.7858
$ws.Range("L74").Value2 = 83335600
$ws.Range("M74").Value2 = -827.7858000000001
$ws.Range("N74").Value2 = -83337348
$ws.Range("H77").Value2 = 25001872
$ws.Range("I77").Value2 = 1701.7858
$ws.Range("J77").Value2 = 83335600
$ws.Range("K77").Value2 = 8508.929
$ws.Range("L77").Value2 = 416678000
$ws.Range("M77").Value2 = -4140.929
$ws.Range("N77").Value2 = -416686736
$ws.Range("H116").Value2 = 1240.909
$ws.Range("I116").Value2 = 801.5714
$ws.Range("K116").Value2 = 801.5714
$ws.Range("M116").Value2 = 1492.4286
$ws.Range("H132").Value2 = 2176612.8
$ws.Range("I132").Value2 = 1977.5555
$ws.Range("J132").Value2 = 10005300
$ws.Range("K132").Value2 = 5932.666499999999
$ws.Range("L132").Value2 = 30015900
$ws.Range("M132").Value2 = -3402.666499999999
$ws.Range("N132").Value2 = -30020960
$ws.Range("H136").Value2 = 4477.4414
$ws.Range("I136").Value2 = 4636.5483
$ws.Range("J136").Value2 = 2833.3333
$ws.Range("K136").Value2 = 13909.6449
$ws.Range("L136").Value2 = 8499.999899999999
$ws.Range("M136").Value2 = -11359.6449
$ws.Range("N136").Value2 = -13599.9999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1240.909
$ws.Range("I3").Value2 = 801.5714
$ws.Range("K3").Value2 = 801.5714
$ws.Range("M3").Value2 = -687.5714
$ws.Range("H12").Value2 = 296.66666
$ws.Range("I12").Value2 = 296.66666
$ws.Range("K12").Value2 = 296.66666
$ws.Range("M12").Value2 = -128.66666
$ws.Range("H134").Value2 = 5674.9033
$ws.Range("I134").Value2 = 7347.1
$ws.Range("J134").Value2 = 2634.5454
$ws.Range("K134").Value2 = 22041.3
$ws.Range("L134").Value2 = 7903.6362
$ws.Range("M134").Value2 = -19506.3
$ws.Range("N134").Value2 = -12973.6362

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 220383.72
$ws.Range("I31").Value2 = 1738.4147
$ws.Range("J31").Value2 = 519198.97
$ws.Range("K31").Value2 = 1738.4147
$ws.Range("L31").Value2 = 519198.97
$ws.Range("M31").Value2 = -1443.4147
$ws.Range("N31").Value2 = -519788.97
$ws.Range("H34").Value2 = 220383.72
$ws.Range("I34").Value2 = 1738.4147
$ws.Range("J34").Value2 = 519198.97
$ws.Range("K34").Value2 = 1738.4147
$ws.Range("L34").Value2 = 519198.97
$ws.Range("M34").Value2 = -1536.4147
$ws.Range("N34").Value2 = -519602.97
$ws.Range("H58").Value2 = 2733268
$ws.Range("I58").Value2 = 3704288.8
$ws.Range("J58").Value2 = 2272.25
$ws.Range("K58").Value2 = 3704288.8
$ws.Range("L58").Value2 = 2272.25
$ws.Range("M58").Value2 = -3704085.8
$ws.Range("N58").Value2 = -2678.25
$ws.Range("H86").Value2 = 1601.7084
$ws.Range("I86").Value2 = 1472.0588
$ws.Range("J86").Value2 = 1916.5714
$ws.Range("K86").Value2 = 1472.0588
$ws.Range("L86").Value2 = 1916.5714
$ws.Range("M86").Value2 = -349.0588
$ws.Range("N86").Value2 = -4162.5714
$ws.Range("H89").Value2 = 1601.7084
$ws.Range("I89").Value2 = 1472.0588
$ws.Range("J89").Value2 = 1916.5714
$ws.Range("K89").Value2 = 7360.294
$ws.Range("L89").Value2 = 9582.857
$ws.Range("M89").Value2 = -1744.294
$ws.Range("N89").Value2 = -20814.857
$ws.Range("H107").Value2 = 15152258
$ws.Range("I107").Value2 = 22222670
$ws.Range("J107").Value2 = 1374.4286
$ws.Range("K107").Value2 = 22222670
$ws.Range("L107").Value2 = 1374.4286
$ws.Range("M107").Value2 = -22220750
$ws.Range("N107").Value2 = -5214.4286
$ws.Range("H132").Value2 = 2858564.2
$ws.Range("I132").Value2 = 3847413.2
$ws.Range("J132").Value2 = 1889.8334
$ws.Range("K132").Value2 = 11542239.6
$ws.Range("L132").Value2 = 5669.5002
$ws.Range("M132").Value2 = -11539709.6
$ws.Range("N132").Value2 = -10729.5002
$ws.Range("H134").Value2 = 6538002.5
$ws.Range("I134").Value2 = 12823468
$ws.Range("J134").Value2 = 1118.44
$ws.Range("K134").Value2 = 38470404
$ws.Range("L134").Value2 = 3355.32
$ws.Range("M134").Value2 = -38467869
$ws.Range("N134").Value2 = -8425.32
$ws.Range("H136").Value2 = 2733268
$ws.Range("I136").Value2 = 3704288.8
$ws.Range("J136").Value2 = 2272.25
$ws.Range("K136").Value2 = 11112866.4
$ws.Range("L136").Value2 = 6816.75
$ws.Range("M136").Value2 = -11110316.4
$ws.Range("N136").Value2 = -11916.75
$ws.Range("H141").Value2 = 33348.11
$ws.Range("J141").Value2 = 33348.11
$ws.Range("L141").Value2 = 33348.11
$ws.Range("N141").Value2 = -43708.11

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value2 = 4761988
$ws.Range("J23").Value2 = 96.2
$ws.Range("L23").Value2 = 288.6
$ws.Range("N23").Value2 = -758.6
$ws.Range("H29").Value2 = 887.6923
$ws.Range("I29").Value2 = 110
$ws.Range("J29").Value2 = 1554.2858
$ws.Range("K29").Value2 = 330
$ws.Range("L29").Value2 = 4662.857400000001
$ws.Range("M29").Value2 = -53
$ws.Range("N29").Value2 = -5216.857400000001
$ws.Range("H80").Value2 = 190
$ws.Range("I80").Value2 = 190
$ws.Range("J80").Value2 = 0
$ws.Range("K80").Value2 = 570
$ws.Range("L80").Value2 = 0
$ws.Range("M80").Value2 = 366
$ws.Range("H83").Value2 = 190
$ws.Range("I83").Value2 = 190
$ws.Range("J83").Value2 = 0
$ws.Range("K83").Value2 = 1710
$ws.Range("L83").Value2 = 0
$ws.Range("M83").Value2 = 2970
$ws.Range("N80").ClearContents()
$ws.Range("N83").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value2 = 5913334
$ws.Range("I7").Value2 = 17500000
$ws.Range("J7").Value2 = 120000.75
$ws.Range("K7").Value2 = 17500000
$ws.Range("L7").Value2 = 120000.75
$ws.Range("M7").Value2 = -17499888
$ws.Range("N7").Value2 = -120224.75
$ws.Range("H8").Value2 = 5913334
$ws.Range("I8").Value2 = 17500000
$ws.Range("J8").Value2 = 120000.75
$ws.Range("K8").Value2 = 17500000
$ws.Range("L8").Value2 = 120000.75
$ws.Range("M8").Value2 = -17499861
$ws.Range("N8").Value2 = -120278.75
$ws.Range("H126").Value2 = 5787.6206
$ws.Range("I126").Value2 = 12845.777
$ws.Range("J126").Value2 = 2611.45
$ws.Range("K126").Value2 = 38537.331
$ws.Range("L126").Value2 = 7834.349999999999
$ws.Range("M126").Value2 = -36067.331
$ws.Range("N126").Value2 = -12774.35
$ws.Range("H132").Value2 = 9806593
$ws.Range("I132").Value2 = 11113751
$ws.Range("J132").Value2 = 2907
$ws.Range("K132").Value2 = 33341253
$ws.Range("L132").Value2 = 8721
$ws.Range("M132").Value2 = -33338723
$ws.Range("N132").Value2 = -13781

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value2 = 55557412
$ws.Range("I68").Value2 = 1958.3334
$ws.Range("J68").Value2 = 166668320
$ws.Range("K68").Value2 = 1958.3334
$ws.Range("L68").Value2 = 166668320
$ws.Range("M68").Value2 = -1209.3334
$ws.Range("N68").Value2 = -166669818
$ws.Range("H71").Value2 = 55557412
$ws.Range("I71").Value2 = 1958.3334
$ws.Range("J71").Value2 = 166668320
$ws.Range("K71").Value2 = 9791.666999999999
$ws.Range("L71").Value2 = 833341600
$ws.Range("M71").Value2 = -6047.666999999999
$ws.Range("N71").Value2 = -833349088
$ws.Range("H132").Value2 = 10279104
$ws.Range("I132").Value2 = 15270921
$ws.Range("J132").Value2 = 1833.6471
$ws.Range("K132").Value2 = 45812763
$ws.Range("L132").Value2 = 5500.9413
$ws.Range("M132").Value2 = -45810233
$ws.Range("N132").Value2 = -10560.9413
$ws.Range("H136").Value2 = 8444.736999999999
$ws.Range("I136").Value2 = 6690.231
$ws.Range("J136").Value2 = 12246.167
$ws.Range("K136").Value2 = 20070.693
$ws.Range("L136").Value2 = 36738.501
$ws.Range("M136").Value2 = -17520.693
$ws.Range("N136").Value2 = -41838.501

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 3167
$ws.Range("I62").Value2 = 3167
$ws.Range("J62").Value2 = 0
$ws.Range("K62").Value2 = 3167
$ws.Range("L62").Value2 = 0
$ws.Range("M62").Value2 = -2543
$ws.Range("H65").Value2 = 3167
$ws.Range("I65").Value2 = 3167
$ws.Range("J65").Value2 = 0
$ws.Range("K65").Value2 = 15835
$ws.Range("L65").Value2 = 0
$ws.Range("M65").Value2 = -12715
$ws.Range("H132").Value2 = 790.1539
$ws.Range("I132").Value2 = 579.8837
$ws.Range("J132").Value2 = 1201.1364
$ws.Range("K132").Value2 = 1739.6511
$ws.Range("L132").Value2 = 3603.4092
$ws.Range("M132").Value2 = 790.3489
$ws.Range("N132").Value2 = -8663.4092
$ws.Range("H136").Value2 = 12580929
$ws.Range("I136").Value2 = 2382.9778
$ws.Range("J136").Value2 = 83335250
$ws.Range("K136").Value2 = 7148.9334
$ws.Range("L136").Value2 = 250005750
$ws.Range("M136").Value2 = -4598.9334
$ws.Range("N136").Value2 = -250010850
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()
